# Weekly "Ajo" (garlic) price update: insert a new weekly record row before
# the existing row 140, shifting the following rows (old 140-143) down to
# (141-144). This matches a new "Feria Lagunitas de Puerto Montt" entry for
# the week of 2021-09-09.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 140; existing rows 140-143
# shift down to 141-144, carrying all of their original values/format.
$ws.Rows.Item(140).Insert()

# Populate the newly inserted row 140 with the new weekly record.
$ws.Cells.Item(140, 1).Value = 4
$ws.Cells.Item(140, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(140, 3).Value = "Los Lagos"
$ws.Cells.Item(140, 4).Value = 44448
$ws.Cells.Item(140, 5).Value = 10
$ws.Cells.Item(140, 6).Value = 100112003
$ws.Cells.Item(140, 7).Value = "Ajo"
$ws.Cells.Item(140, 8).Value = "Chino"
$ws.Cells.Item(140, 9).Value = "Primera"
$ws.Cells.Item(140, 10).Value = 120
$ws.Cells.Item(140, 11).Value = 18000
$ws.Cells.Item(140, 12).Value = 18000
$ws.Cells.Item(140, 13).Value = 18000
$ws.Cells.Item(140, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(140, 15).Value = "China"
$ws.Cells.Item(140, 16).Value = 1800
$ws.Cells.Item(140, 17).Value = 10
$ws.Cells.Item(140, 18).Value = "Hortaliza"
